# edit.ps1
# Applies the GitHub Actions "Updated cryptos list" data refresh to the
# cryptos worksheet: updates Price (column D) and Volume(1h) (column E)
# for every data row, and for the three rows whose coin ranking changed
# also updates Coin (column B) and Link (column C) so the row order keeps
# matching the new ranking:
#   - rows 14-16: Chainlink / WrappedEther / BinanceUSD re-ordered
#   - rows 38-39: FraxShare / TrustWalletToken re-ordered
#
# The Price column stores values such as "1.024" or "28.556.37" as plain
# text (they are not real numbers - some even contain two '.' separators).
# Excel auto-detects a lone-dot numeric-looking string and silently stores
# it as a Number, so those values are written with a leading apostrophe
# (forces a text entry, exactly like a user typing '1.024 into the cell)
# and the cell style is immediately reset back to Normal so no stray
# number-format / quote-prefix formatting is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row=2; B=$null; C=$null; D='28.556.37'; E='  +2.06%  ' },
    @{ Row=3; B=$null; C=$null; D='1.878.53'; E='  +0.72%  ' },
    @{ Row=4; B=$null; C=$null; D='1.024'; E='  +2.27%  ' },
    @{ Row=5; B=$null; C=$null; D='318.61'; E='  +2.21%  ' },
    @{ Row=6; B=$null; C=$null; D='1.024'; E='  +2.48%  ' },
    @{ Row=7; B=$null; C=$null; D='0.5157'; E='  +1.29%  ' },
    @{ Row=8; B=$null; C=$null; D='0.3975'; E='  +2.70%  ' },
    @{ Row=9; B=$null; C=$null; D='0.08375'; E='  +0.39%  ' },
    @{ Row=10; B=$null; C=$null; D='1.115'; E='  +0.43%  ' },
    @{ Row=11; B=$null; C=$null; D='42.21'; E='  +2.27%  ' },
    @{ Row=12; B=$null; C=$null; D='6.270'; E='  +1.17%  ' },
    @{ Row=13; B=$null; C=$null; D='20.58'; E='  +0.23%  ' },
    @{ Row=14; B='WrappedEther'; C='https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; D='1.822.24'; E='  -1.58%  ' },
    @{ Row=15; B='BinanceUSD'; C='https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'; D='1.027'; E='  +2.70%  ' },
    @{ Row=16; B='Chainlink'; C='https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; D='7.250'; E='  +0.62%  ' },
    @{ Row=17; B=$null; C=$null; D='0.00001112'; E='  +1.56%  ' },
    @{ Row=18; B=$null; C=$null; D='91.53'; E='  +1.25%  ' },
    @{ Row=19; B=$null; C=$null; D='0.06786'; E='  +2.01%  ' },
    @{ Row=20; B=$null; C=$null; D='17.78'; E='  +0.76%  ' },
    @{ Row=21; B=$null; C=$null; D='1.024'; E='  +2.47%  ' },
    @{ Row=22; B=$null; C=$null; D='5.981'; E='  +0.08%  ' },
    @{ Row=23; B=$null; C=$null; D='28.594.55'; E='  +2.03%  ' },
    @{ Row=24; B=$null; C=$null; D='11.20'; E='  +1.34%  ' },
    @{ Row=25; B=$null; C=$null; D='2.294'; E='  +2.65%  ' },
    @{ Row=26; B=$null; C=$null; D='162.70'; E='  +2.68%  ' },
    @{ Row=27; B=$null; C=$null; D='2.044.67'; E='  -1.80%  ' },
    @{ Row=28; B=$null; C=$null; D='20.77'; E='  +1.43%  ' },
    @{ Row=29; B=$null; C=$null; D='2.365'; E='  -3.57%  ' },
    @{ Row=30; B=$null; C=$null; D='127.75'; E='  +2.70%  ' },
    @{ Row=31; B=$null; C=$null; D='0.1048'; E='  -0.42%  ' },
    @{ Row=32; B=$null; C=$null; D='1.043'; E='  +1.54%  ' },
    @{ Row=33; B=$null; C=$null; D='5.837'; E='  +0.99%  ' },
    @{ Row=34; B=$null; C=$null; D='3.645'; E='  +1.59%  ' },
    @{ Row=35; B=$null; C=$null; D='0.02438'; E='  +0.19%  ' },
    @{ Row=36; B=$null; C=$null; D='0.06514'; E='  +0.03%  ' },
    @{ Row=37; B=$null; C=$null; D='0.2193'; E='  +0.66%  ' },
    @{ Row=38; B='TrustWalletToken'; C='https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; D='1.284'; E='  +5.79%  ' },
    @{ Row=39; B='FraxShare'; C='https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; D='8.918'; E='  -5.95%  ' },
    @{ Row=40; B=$null; C=$null; D='1.190'; E='  +0.21%  ' },
    @{ Row=41; B=$null; C=$null; D='0.6460'; E='  +0.37%  ' },
    @{ Row=42; B=$null; C=$null; D='5.043'; E='  +1.90%  ' },
    @{ Row=43; B=$null; C=$null; D='11.28'; E='  +0.55%  ' },
    @{ Row=44; B=$null; C=$null; D='0.6048'; E='  -0.45%  ' },
    @{ Row=45; B=$null; C=$null; D='13.19'; E='  +1.72%  ' },
    @{ Row=46; B=$null; C=$null; D='3.737'; E='  +2.32%  ' },
    @{ Row=47; B=$null; C=$null; D='1.232'; E='  -3.45%  ' },
    @{ Row=48; B=$null; C=$null; D='1.999'; E='  +0.43%  ' },
    @{ Row=49; B=$null; C=$null; D='1.211'; E='  -0.76%  ' },
    @{ Row=50; B=$null; C=$null; D='122.23'; E='  +1.90%  ' },
    @{ Row=51; B=$null; C=$null; D='0.06870'; E='  -0.06%  ' }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.B) { $ws.Cells.Item($r, 2).Value = $u.B }
    if ($null -ne $u.C) { $ws.Cells.Item($r, 3).Value = $u.C }
    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($r, 4)
        $cell.Value = "'" + $u.D
        $cell.Style = "Normal"
    }
    if ($null -ne $u.E) { $ws.Cells.Item($r, 5).Value = $u.E }
}
